# ============================================================================
# Apply DCF workbook edit: update Assumptions sheet, add a new "Results"
# sheet that surfaces the key DCF outputs, and re-point helper formulas at
# the DCF sheet's actual discounted-cash-flow column (S) instead of the old
# placeholder C19:H19 / H19 references.
# ============================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Assumptions sheet tweaks
# ---------------------------------------------------------------------------
$assump = $wb.Worksheets("Assumptions")

# Terminal discount rate: 16% -> 14%
$assump.Range("B7").Value = 0.14

# Relabel + repoint the summary rows to the real DCF columns / new key names
$assump.Range("A9").Value = "discounted_window_cash_flow"
$assump.Range("B9").Formula = "=SUM(DCF!S3:S8)"

$assump.Range("A10").Value = "discounted_terminal_cash_flow"
$assump.Range("B10").Formula = "=DCF!S8*(1+B6)/(B7-B6)"

$assump.Range("A11").Value = "net_debt"

$assump.Range("A12").Value = "equity_value"

# New rows: fully diluted share count + per-share equity value
$assump.Range("A13").Value = "fully_diluted_shares"
$assump.Range("B13").Value = 10

$assump.Range("A14").Value = "equity_value_per_share"
$assump.Range("B14").Formula = "=B12/B13"
$assump.Range("B14").NumberFormat = $assump.Range("B12").NumberFormat

# Widen column A so the longer key names fit
$assump.Columns.Item(1).ColumnWidth = 29

# ---------------------------------------------------------------------------
# 2. New "Results" sheet, placed after "Assumptions"
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$results = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$results.Name = "Results"

$results.Range("A1").Value = "Key"
$results.Range("B1").Value = "Value"

$currencyFormat = $assump.Range("B9").NumberFormat

function Set-ResultRow($row, $key, $formulaOrValue, $isFormula, $valueStyle) {
    $keyCell = $results.Range("A$row")
    $keyCell.Value = $key
    $keyCell.Font.Color = 0

    $valCell = $results.Range("B$row")
    if ($isFormula) {
        $valCell.Formula = $formulaOrValue
    } else {
        $valCell.Value = $formulaOrValue
    }
    $valCell.Font.Color = 0
    if ($valueStyle -eq "currency") {
        $valCell.NumberFormat = $currencyFormat
    }
}

Set-ResultRow 2 "discounted_window_cash_flow" "=SUM(DCF!S3:S8)" $true "currency"
Set-ResultRow 3 "discounted_terminal_cash_flow" "=DCF!S8*(1+Assumptions!`$B`$6)/(Assumptions!`$B`$7-Assumptions!`$B`$6)/(1+Assumptions!`$B`$3)^DCF!`$A`$8" $true "currency"
Set-ResultRow 4 "net_debt" 2 $false "currency"
Set-ResultRow 5 "enterprise_value" "=B2+B3" $true "currency"
Set-ResultRow 6 "equity_value" "=B5-B4" $true "currency"
Set-ResultRow 7 "fully_diluted_shares" 10 $false "plain"
Set-ResultRow 8 "equity_value_per_share" "=B6/B7" $true "currency"

$results.Columns.Item(1).ColumnWidth = 27

# ---------------------------------------------------------------------------
# 3. Selections / active sheet (cosmetic, matches author's final view state)
# ---------------------------------------------------------------------------
$dcf = $wb.Worksheets("DCF")
$dcf.Range("G23").Select()

$assump.Range("B8").Select()

$results.Range("F15").Select()
$results.Activate()
